$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new IR Sensor row content on row 38 (the existing "Servo" row), per the
# author's attempt to wire up a rear IR obstacle sensor.
$ws.Range("A38").Value = "IR Sensor"
$ws.Range("C38").Value = "Sends light signals to to detect obstacles. Used for rear obstacle detection"

# Widen column C to fit the new, longer description text (40.66 -> 60.66 chars).
# The host's ColumnWidth setter quantizes through a pixel round-trip, so feed
# it a value whose rounded pixel width lands back on the target character
# width as closely as possible.
$ws.Columns.Item(3).ColumnWidth = 59.83
